# The deck's Design was changed from the custom "Integral" (Red Violet)
# theme to the built-in default "Office Theme" - i.e. the user picked a
# different theme from the Design gallery. In OOXML terms this means the
# colour scheme that is bound to the (single) slide master - and hence to
# every slide - switches from the "Red Violet" palette to the standard
# "Office" palette.
#
# PowerPoint exposes the twelve slots of a theme's colour scheme through
# Slide.ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink,
# in that order). Writing to ThemeColorScheme.Item(n).RGB updates the
# presentation's theme part in place, so re-pointing every slot at the
# standard Office colours reproduces the effect of applying the Office
# Theme design to the whole deck.

$p = $ppt.ActivePresentation

# Standard "Office" theme colour scheme (RRGGBB -> COM RGB long = R + G*256 + B*65536)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slideCount = $p.Slides.Count
for ($s = 1; $s -le $slideCount; $s++) {
    $slide = $p.Slides.Item($s)
    $tcs = $slide.ThemeColorScheme
    for ($i = 1; $i -le $tcs.Count; $i++) {
        $tcs.Item($i).RGB = $officeColors[$i - 1]
    }
}
